$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 0: swap header labels in row 1 (A1 <-> B1) ---
$a1 = $ws.Range("A1").Value()
$b1 = $ws.Range("B1").Value()
$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

# --- Step 1: move column B (API description) content into column A for rows 2-25, preserving Bs original style ---
for ($r = 2; $r -le 25; $r++) {
  $bCell = $ws.Cells.Item($r, 2)
  $aCell = $ws.Cells.Item($r, 1)
  $val = $bCell.Value()
  $bCell.Copy() | Out-Null
  $aCell.PasteSpecial(-4122) | Out-Null
  $aCell.Value = $val
}
$excel.CutCopyMode = $false

# --- Step 2: restyle B cells that must become style "s1" (short-description style) before setting new text ---
# Use B4 (already s1) as the format donor.
$donor = $ws.Cells.Item(4, 2)
$donor.Copy() | Out-Null
$ws.Cells.Item(3, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(5, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(6, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(7, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(13, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(15, 2).PasteSpecial(-4122) | Out-Null
$donor.Copy() | Out-Null
$ws.Cells.Item(16, 2).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 3: set new short-description text (or clear) for column B, rows 2-25 ---
$ws.Cells.Item(2, 2).ClearContents() | Out-Null
$ws.Cells.Item(3, 2).Value = "See my after hours energy consumption sorted by end use"
$ws.Cells.Item(4, 2).Value = "Compare my after hours energy consumption to my total energy consumption"
$ws.Cells.Item(5, 2).Value = "See my total energy consumption"
$ws.Cells.Item(6, 2).Value = "See my total energy consumption and end use breakdown"
$ws.Cells.Item(7, 2).Value = "See my total energy consumption for a specific end use"
$ws.Cells.Item(8, 2).Value = "Tell me how much energy was used by a specific piece of equipment "
$ws.Cells.Item(9, 2).ClearContents() | Out-Null
$ws.Cells.Item(10, 2).ClearContents() | Out-Null
$ws.Cells.Item(11, 2).ClearContents() | Out-Null
$ws.Cells.Item(12, 2).Value = "Compare my after hours energy consumption to my total energy consumption"
$ws.Cells.Item(13, 2).Value = "Compare my after hours energy consumption to my total energy consumption"
$ws.Cells.Item(14, 2).Value = "Identify the maximum power draw and time it happened for a piece of equipment"
$ws.Cells.Item(15, 2).Value = "Identify the maximum energy consumption and the approximate time it occured in my building"
$ws.Cells.Item(16, 2).Value = "See my monthly energy consumption (see, look at, inspect, compare, analyze: Verb + what they want to do)"
$ws.Cells.Item(17, 2).ClearContents() | Out-Null
$ws.Cells.Item(18, 2).ClearContents() | Out-Null
$ws.Cells.Item(19, 2).ClearContents() | Out-Null
$ws.Cells.Item(20, 2).ClearContents() | Out-Null
$ws.Cells.Item(21, 2).ClearContents() | Out-Null
$ws.Cells.Item(22, 2).ClearContents() | Out-Null
$ws.Cells.Item(23, 2).ClearContents() | Out-Null
$ws.Cells.Item(24, 2).ClearContents() | Out-Null
$ws.Cells.Item(25, 2).ClearContents() | Out-Null

# --- Step 4: column widths ---
$ws.Columns.Item(1).ColumnWidth = 110.83333333333333
$ws.Columns.Item(2).ColumnWidth = 55.166666666666664

# --- Step 5: selection ---
$ws.Range("A6").Select() | Out-Null

Write-Output "done"